# Auto-generated Excel COM-interop script
# Applies the cell-value updates to the cryptos worksheet described by the commit diff
# (live crypto price/volume refresh + a few re-ordered rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.237.23'
$ws.Range("D3").Value = '3.870.54'
$ws.Range("E3").Value = '  +3.21%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'604.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").Value = "'165.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.04%  '
$ws.Range("D7").Value = '3.868.12'
$ws.Range("E7").Value = '  +3.22%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = "'0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.24%  '
$ws.Range("E10").Value = '  +0.98%  '
$ws.Range("D11").Value = "'6.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("E13").Value = '  -2.46%  '
$ws.Range("D14").Value = "'0.0000247"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.74%  '
$ws.Range("D15").Value = '4.517.76'
$ws.Range("E15").Value = '  +3.22%  '
$ws.Range("D16").Value = '3.870.82'
$ws.Range("E16").Value = '  +3.28%  '
$ws.Range("D17").Value = '69.400.84'
$ws.Range("E17").Value = '  +0.38%  '
$ws.Range("E18").Value = '  +3.96%  '
$ws.Range("D19").Value = "'11.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.48%  '
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("D21").Value = "'17.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").Value = "'490.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = "'0.0000160"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.71%  '
$ws.Range("D25").Value = "'84.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("E26").Value = '  -1.56%  '
$ws.Range("E27").Value = '  -0.55%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = "'10.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.63%  '
$ws.Range("D30").Value = "'3.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = "'2.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.06%  '
$ws.Range("B33").Value = 'WrappedeETH'
$ws.Range("C33").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D33").Value = '4.023.21'
$ws.Range("E33").Value = '  +3.31%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = "'32.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.90%  '
$ws.Range("D35").Value = '3.816.31'
$ws.Range("E35").Value = '  +3.62%  '
$ws.Range("E36").Value = '  -1.09%  '
$ws.Range("E37").Value = '  +2.27%  '
$ws.Range("E38").Value = '  +4.28%  '
$ws.Range("E39").Value = '  +0.80%  '
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("E41").Value = '  -1.38%  '
$ws.Range("E42").Value = '  +1.48%  '
$ws.Range("D43").Value = "'441.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.83%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'27.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +19.15%  '
$ws.Range("B47").Value = 'Cosmos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D47").Value = "'8.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.90%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = "'1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("D49").Value = '2.862.29'
$ws.Range("E49").Value = '  +2.40%  '
$ws.Range("D50").Value = "'143.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.37%  '
$ws.Range("E51").Value = '  +1.37%  '

Write-Host "Applied crypto list update."
